# Daily attendance processing - move "System" entry from the front to the
# back of the comma-separated "Recorded By" list in column G.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$rowCount = $usedRange.Rows.Count()

for ($i = 2; $i -le $rowCount; $i++) {
    $cell = $ws.Cells.Item($i, 7)
    $val = $cell.Value()
    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ", "
        if ($parts.Count -gt 1 -and $parts[0] -eq "System") {
            $rest = $parts[1..($parts.Count - 1)]
            $newParts = $rest + $parts[0]
            $newVal = $newParts -join ", "
            $cell.Value = $newVal
        }
    }
}
